# simplify steel description (remove RME)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B2 holds the multi-line "industrial" mapping description.
# Remove the "/RME" token from the steel (S) line:
#   "3% S/LFM+CDL/RME/H:1" -> "3% S/LFM+CDL/H:1"
$cell = $ws.Range("B2")
$text = $cell.Value2
$newText = $text.Replace("/RME", "")
$cell.Value = $newText

# Move the active selection to B2 (previously B3).
$ws.Range("B2").Select()
